# Test data file 6th jan
# data file after flushing the database
#
# Refresh the Sheet1 "TestData" values to the new (post-flush) test
# dataset, drop the now-unused Hyperlink cell style, and move the
# sheet view/selection back to the top-left area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1-2: login creds -------------------------------------------------
$ws.Range("B1").Value = "ruchira.m"
$ws.Range("B2").Value = "ruchira.m,1234"

# --- Row 4: CreateInvoice data row ----------------------------------------
$ws.Range("B4").Value = "Amanora"
$ws.Range("D4").Value = "bangalore"
$ws.Range("H4").Value = "test product"
$ws.Range("I4").Value = "Department1"
$ws.Range("J4").Value = "Fuel"
$ws.Range("K4").Value = "inv desc2"
$ws.Range("L4").Value = "measure1"
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 10.44444
$ws.Range("O4").Value = "message1"
$ws.Range("P4").Value = "memo1"
$ws.Range("Q4").Value = "NetChain2-AP New: Invoice"

# --- Row 7: CreateVendor data row -----------------------------------------
$ws.Range("E7").Value = "Vendor 1"
$ws.Range("F7").Value = "abcd@gmail.com"

# F7 used to look like a hyperlink (underline + themed colour) purely via
# the "Hyperlink" cell style, with no actual HYPERLINK()/hyperlink
# relationship behind it. Clear it back to Normal and drop the now-unused
# named style from the workbook.
$ws.Range("F7").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

# --- sheet view: selection moves to A6, no forced horizontal scroll -------
$ws.Activate() | Out-Null
$ws.Range("A6").Select() | Out-Null
